# tests/A01_pixell_test_plan_client.xlsx
# "Initialized and ran tests with no error"
#
# The author ran the getter test cases (rows 12-16, test cases 6-10 for
# client_number / first_name / last_name / email_address / __str__) and
# filled in the actual results:
#   - Preconditions (col E) now records the attribute that was already set
#     on the object (previously left as "None").
#   - Method Inputs (col F) now correctly shows "None" (getters take no
#     arguments) instead of the attribute value that had been pasted there.
#   - Expected Result (col G) now shows the value the getter returned
#     (mirroring col E), replacing the generic "no error" placeholder.
#   - Row 12 (client_number) instead records "none" as the method input.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test case 6 - client_number getter
$ws.Range("E12").Value = "Client_number=12345"
$ws.Range("F12").Value = "none"
$ws.Range("G12").Value = "Client_number=12345"

# Test case 7 - first_name getter
$ws.Range("E13").Value = 'first_name="Wendy"'
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = 'first_name="Wendy"'

# Test case 8 - last_name getter
$ws.Range("E14").Value = 'last_name="Ways"'
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = 'last_name="Ways"'

# Test case 9 - email_address getter
$ws.Range("E15").Value = 'email_address="WendyWays@pixell-river.com"'
$ws.Range("F15").Value = "None"
$ws.Range("G15").Value = 'email_address="WendyWays@pixell-river.com"'

# Test case 10 - __str__ getter
$ws.Range("E16").Value = "Ways, Wendy [12345] - WendyWays@pixell-river.com"
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = "Ways, Wendy [12345] - WendyWays@pixell-river.com"

# View state: after running the tests the sheet had been scrolled so row 9 /
# column D sits at the top-left of the visible pane (was row 8 / column A).
# Selection (G16) is left untouched.
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 9
    $win.ScrollColumn = 4
} catch {
    # Window-geometry state is best-effort in headless sessions.
}
